# daily auto push: 2026-01-26 18:52 UTC
# The log sheet is a rolling window of timestamped entries. Two new rows
# are inserted at the top of the data block (row 703) with the latest
# entries, which pushes the existing rows 703:744 down to 705:746.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right before the existing "2026/12/29" block so
# everything currently at rows 703:744 slides down to 705:746.
$ws.Rows.Item(703).Insert()
$ws.Rows.Item(703).Insert()

# New row 703: 2026/01/26 (Mon)
$ws.Range("A703").Value = "'2026/01/26"
$ws.Range("A703").Style = "Normal"
$ws.Range("B703").Value = "月"
$ws.Range("C703").Value = 22
$ws.Range("D703").Value = 22

# New row 704: 2026/01/27 (Mon)
$ws.Range("A704").Value = "'2026/01/27"
$ws.Range("A704").Style = "Normal"
$ws.Range("B704").Value = "月"
$ws.Range("C704").Value = 1
$ws.Range("D704").Value = 23
